$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.886.35'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.417.23'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '551.44'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.18'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +3.56%  '
$ws.Range('E9').Value = '  -1.94%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.70'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.147'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  -2.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.54'
$ws.Range('E13').Value = '  +3.70%  '
$ws.Range('D14').Value = '2.845.32'
$ws.Range('D15').Value = '59.826.36'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '2.408.43'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.31'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '329.04'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.64'
$ws.Range('E21').Value = '  -4.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.54'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.62'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '0.0₃0775'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.77'
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.39'
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.11'
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.61'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.01'
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.22'
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.409'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '313.62'
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '138.64'
$ws.Range('E42').Value = '  -2.68%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '19.53'
$ws.Range('E45').Value = '  +2.51%  '
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('B48').Value = 'Polygon'
$ws.Range('C48').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.388'
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '17.69'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.06'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('E51').Value = '  -1.46%  '
